$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = '62.957.07'
$ws.Cells.Item(2, 5).Value = '  +3.86%  '
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = '2.468.37'
$ws.Cells.Item(3, 5).Value = '  +5.40%  '
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = '0.998'
$ws.Cells.Item(4, 5).Value = '  -0.13%  '
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '563.86'
$ws.Cells.Item(5, 5).Value = '  +2.56%  '
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '142.39'
$ws.Cells.Item(6, 5).Value = '  +8.31%  '
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '0.999'
$ws.Cells.Item(7, 5).Value = '  -0.09%  '
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '0.587'
$ws.Cells.Item(8, 5).Value = '  +1.04%  '
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '2.467.33'
$ws.Cells.Item(9, 5).Value = '  +5.48%  '
$ws.Cells.Item(10, 5).Value = '  +2.57%  '
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '5.69'
$ws.Cells.Item(11, 5).Value = '  +1.10%  '
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '0.152'
$ws.Cells.Item(12, 5).Value = '  +1.51%  '
$ws.Cells.Item(13, 5).Value = '  +4.10%  '
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '26.46'
$ws.Cells.Item(14, 5).Value = '  +10.45%  '
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '2.899.28'
$ws.Cells.Item(15, 5).Value = '  +5.09%  '
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '62.794.89'
$ws.Cells.Item(16, 5).Value = '  +3.81%  '
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '0.0000141'
$ws.Cells.Item(17, 5).Value = '  +4.34%  '
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '2.459.70'
$ws.Cells.Item(18, 5).Value = '  +5.54%  '
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '11.22'
$ws.Cells.Item(19, 5).Value = '  +4.80%  '
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '339.36'
$ws.Cells.Item(20, 5).Value = '  +7.83%  '
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '4.25'
$ws.Cells.Item(21, 5).Value = '  +3.27%  '
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '6.80'
$ws.Cells.Item(22, 5).Value = '  +2.82%  '
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '1.00'
$ws.Cells.Item(23, 5).Value = '  +0.05%  '
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '65.61'
$ws.Cells.Item(24, 5).Value = '  +2.05%  '
$ws.Cells.Item(25, 5).Value = '  +1.54%  '
$ws.Cells.Item(26, 5).Value = '  +0.26%  '
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '1.50'
$ws.Cells.Item(27, 5).Value = '  +6.62%  '
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '8.04'
$ws.Cells.Item(28, 5).Value = '  +0.82%  '
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '1.39'
$ws.Cells.Item(29, 5).Value = '  +9.78%  '
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '6.83'
$ws.Cells.Item(30, 5).Value = '  +11.87%  '
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '1.84'
$ws.Cells.Item(31, 5).Value = '  +5.22%  '
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '0.0₃0795'
$ws.Cells.Item(32, 5).Value = '  +7.98%  '
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '175.66'
$ws.Cells.Item(33, 5).Value = '  +2.80%  '
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '1.51'
$ws.Cells.Item(34, 5).Value = '  +9.85%  '
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '0.395'
$ws.Cells.Item(35, 5).Value = '  +2.29%  '
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '18.77'
$ws.Cells.Item(36, 5).Value = '  +3.80%  '
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '367.16'
$ws.Cells.Item(37, 5).Value = '  +11.35%  '
$ws.Cells.Item(38, 5).Value = '  +0.00%  '
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '4.38'
$ws.Cells.Item(39, 5).Value = '  +5.72%  '
$ws.Cells.Item(40, 5).Value = '  -0.16%  '
$ws.Cells.Item(41, 5).Value = '  +10.00%  '
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '40.58'
$ws.Cells.Item(42, 5).Value = '  +6.53%  '
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '149.66'
$ws.Cells.Item(43, 5).Value = '  +8.33%  '
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '3.70'
$ws.Cells.Item(44, 5).Value = '  +4.87%  '
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '20.42'
$ws.Cells.Item(45, 5).Value = '  +5.65%  '
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '0.597'
$ws.Cells.Item(46, 5).Value = '  +4.62%  '
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '0.0957'
$ws.Cells.Item(47, 5).Value = '  +0.60%  '
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '0.0515'
$ws.Cells.Item(48, 5).Value = '  +3.10%  '
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '0.0₆0233'
$ws.Cells.Item(49, 5).Value = '  +4.61%  '
$ws.Cells.Item(50, 5).Value = '  +4.07%  '
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '17.85'
$ws.Cells.Item(51, 5).Value = '  +3.86%  '
